$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "合富中国"
$ws.Cells.Item(2, 2).Value = "合富中国"
$ws.Cells.Item(2, 3).Value = "九牧王"
$ws.Cells.Item(3, 1).Value = "航天发展"
$ws.Cells.Item(3, 2).Value = "航天发展"
$ws.Cells.Item(4, 1).Value = "九牧王"
$ws.Cells.Item(4, 2).Value = "九牧王"
$ws.Cells.Item(5, 1).Value = "华胜天成"
$ws.Cells.Item(5, 2).Value = "华胜天成"
$ws.Cells.Item(5, 3).Value = "平潭发展"
$ws.Cells.Item(6, 2).Value = "华夏幸福"
$ws.Cells.Item(6, 3).Value = "多氟多"
$ws.Cells.Item(7, 1).Value = "海南海药"
$ws.Cells.Item(7, 2).Value = "大为股份"
$ws.Cells.Item(8, 1).Value = "榕基软件"
$ws.Cells.Item(8, 2).Value = "中国银行"
$ws.Cells.Item(8, 3).Value = "华夏幸福"
$ws.Cells.Item(9, 2).Value = "中水渔业"
$ws.Cells.Item(9, 3).Value = "榕基软件"
$ws.Cells.Item(10, 1).Value = "华夏幸福"
$ws.Cells.Item(10, 2).Value = "海南海药"
$ws.Cells.Item(10, 3).Value = "大为股份"
$ws.Cells.Item(11, 1).Value = "平潭发展"
$ws.Cells.Item(11, 2).Value = "海马汽车"
$ws.Cells.Item(11, 3).Value = "智能自控"
$ws.Cells.Item(12, 1).Value = "多氟多"
$ws.Cells.Item(12, 2).Value = "浪潮软件"
$ws.Cells.Item(12, 3).Value = "浪潮软件"
$ws.Cells.Item(13, 1).Value = "海马汽车"
$ws.Cells.Item(13, 2).Value = "多氟多"
$ws.Cells.Item(13, 3).Value = "海马汽车"
$ws.Cells.Item(14, 1).Value = "浪潮软件"
$ws.Cells.Item(14, 2).Value = "榕基软件"
$ws.Cells.Item(14, 3).Value = "华映科技"
$ws.Cells.Item(15, 1).Value = "华映科技"
$ws.Cells.Item(15, 2).Value = "华映科技"
$ws.Cells.Item(15, 3).Value = "中水渔业"
$ws.Cells.Item(16, 1).Value = "国联水产"
$ws.Cells.Item(16, 2).Value = "国联水产"
$ws.Cells.Item(16, 3).Value = "海南海药"
$ws.Cells.Item(17, 1).Value = "中国银行"
$ws.Cells.Item(17, 2).Value = "大众公用"
$ws.Cells.Item(17, 3).Value = "天际股份"
$ws.Cells.Item(18, 1).Value = "大众公用"
$ws.Cells.Item(18, 2).Value = "方正科技"
$ws.Cells.Item(18, 3).Value = "金圆股份"
$ws.Cells.Item(19, 1).Value = "腾景科技"
$ws.Cells.Item(19, 2).Value = "盛新锂能"
$ws.Cells.Item(19, 3).Value = "雪人集团"
$ws.Cells.Item(20, 1).Value = "蜂助手"
$ws.Cells.Item(20, 2).Value = "和而泰"
$ws.Cells.Item(20, 3).Value = "大众公用"
$ws.Cells.Item(21, 1).Value = "宜通世纪"
$ws.Cells.Item(21, 2).Value = "常山北明"
$ws.Cells.Item(21, 3).Value = "孚日股份"
